$wb = $excel.ActiveWorkbook

# Map of cell -> new value for the "想去人数" (F) column updates.
# These same updates apply identically to both the "展览" sheet and the
# "全部类型" sheet, which mirror each other's data.
$updates = @{
    "F2"  = 1980
    "F4"  = 120
    "F5"  = 41
    "F7"  = 1650
    "F8"  = 24
    "F9"  = 658
    "F13" = 95
    "F18" = 132
    "F19" = 3832
    "F21" = 20
    "F22" = 434
    "F23" = 355
    "F24" = 709
    "F25" = 469
    "F27" = 31
    "F28" = 1631
    "F30" = 156
    "F31" = 4
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
